$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title
#    ("Play Dazzling Diamonds Slot for Free - Review 2021").
#    We build it by duplicating the run-structure of the existing bold
#    paragraph near the end of the document (it already has the
#    "<empty run><bold run>" pattern we need), pasting it after the
#    title, then rewriting its text.
# ---------------------------------------------------------------------

$count = $d.Paragraphs.Count
$pBoldSource = $d.Paragraphs.Item($count - 1)   # "Play Dazzling Diamonds Slot for Free - Review 2021" (bold)
if ($pBoldSource.Range.Text.TrimEnd() -ne "Play Dazzling Diamonds Slot for Free - Review 2021") {
    throw "Unexpected document layout: source bold paragraph not found where expected."
}
$pBoldSource.Range.Copy()

$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"
$metaPara.Range.Paste()

# Replace the bold run's text with "Meta description"
$metaPara = $d.Paragraphs.Item(2)
$boldRange = $d.Range($metaPara.Range.Start, $metaPara.Range.End - 1)
$boldRange.Text = "Meta description"
$boldRange.Bold = 1

# Append the non-bold remainder of the sentence right after it, inside
# the same paragraph (before the paragraph mark).
$metaPara = $d.Paragraphs.Item(2)
$tailStart = $metaPara.Range.End - 1
$tailRange = $d.Range($tailStart, $tailStart)
$tailRange.InsertAfter(": Read our unbiased review of Dazzling Diamonds slot. Find out how to play for free, pros and cons, and experience this classic slot game with a 100,000 euros jackpot.")
$tailRange.Bold = 0

# ---------------------------------------------------------------------
# 2) Remove the now-duplicated bold paragraph
#    ("Play Dazzling Diamonds Slot for Free - Review 2021") that used to
#    sit near the end of the document, right before the italic blurb.
# ---------------------------------------------------------------------

$count = $d.Paragraphs.Count
$pBoldTail = $d.Paragraphs.Item($count - 1)
if ($pBoldTail.Range.Text.TrimEnd() -ne "Play Dazzling Diamonds Slot for Free - Review 2021") {
    throw "Unexpected document layout: trailing bold paragraph not found where expected."
}
$pBoldTail.Range.Delete()

# ---------------------------------------------------------------------
# 3) Replace the text of the trailing italic paragraph with the new
#    image-generation prompt, keeping its italic formatting intact.
# ---------------------------------------------------------------------

$newPrompt = 'Prompt: Design a cartoon-style feature image for the online slot game "Dazzling Diamonds". The image should feature a happy Maya warrior with glasses. Details: - The image should be in a 2D cartoon style with vibrant colors. - The Maya warrior should be male and have a friendly, approachable expression on his face. - The warrior should be wearing a traditional headdress and clothing, but with a modern twist, such as wearing glasses. - In the background, there should be sparkling diamonds and other precious gems to match the theme of the game. - The title "Dazzling Diamonds" should be prominently displayed in the image. - The overall tone should be playful and inviting, encouraging players to try their luck at the game.'

# Replace just the text of the final (italic) paragraph directly, leaving
# its paragraph mark (and therefore its run/paragraph formatting, as well
# as the leading empty run) untouched.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$expectedOld = "Read our unbiased review of Dazzling Diamonds slot. Find out how to play for free, pros and cons, and experience this classic slot game with a 100,000 euros jackpot."
if ($lastPara.Range.Text.TrimEnd() -ne $expectedOld) {
    throw "Unexpected document layout: trailing italic paragraph not found where expected."
}
$textRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$textRange.Text = $newPrompt
